$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1) Rename the data sheet to reflect the new "as of" date.
$ws.Name = "富国500低波医药消费策略指数 20230320"

# 2) Insert a new row at 22 (pushes the old spacer/summary rows down by one,
#    and Excel auto-adjusts the relative formula references in the rows below).
$ws.Rows.Item(22).Insert()

# 3) The freshly inserted row 22 becomes the new "current year" row and
#    should carry the plain format row 21 (2022) currently has, before that
#    format gets overwritten in the next step.
$ws.Range("A21:L21").Copy()
$ws.Range("A22:L22").PasteSpecial(-4122)

# 4) Row 21 (2022) becomes the "previous, now-highlighted" year; give it the
#    highlighted format used historically for such rows (same pattern as
#    row 5, which already carries the 16/17/20 style triple).
$ws.Range("A5:L5").Copy()
$ws.Range("A21:L21").PasteSpecial(-4122)

# 5) Update the 2022 figures (row 21) with the revised numbers.
$ws.Range("A21").Value = 2022
$ws.Range("B21").Value = 3871
$ws.Range("C21").Value = 5864
$ws.Range("D21").Value = 14690
$ws.Range("E21").Value = 13504
$ws.Range("F21").Value = 21573
$ws.Range("H21").Formula = "=(B21-B20)/B20"
$ws.Range("I21").Formula = "=(C21-C20)/C20"
$ws.Range("J21").Formula = "=(D21-D20)/D20"
$ws.Range("K21").Formula = "=(E21-E20)/E20"
$ws.Range("L21").Formula = "=(F21-F20)/F20"

# 6) Populate the new 2023 row (row 22).
$ws.Range("A22").Value = 2023
$ws.Range("B22").Value = 3939
$ws.Range("C22").Value = 6160
$ws.Range("D22").Value = 16070
$ws.Range("E22").Value = 13391
$ws.Range("F22").Value = 21034
$ws.Range("H22").Formula = "=(B22-B21)/B21"
$ws.Range("I22").Formula = "=(C22-C21)/C21"
$ws.Range("J22").Formula = "=(D22-D21)/D21"
$ws.Range("K22").Formula = "=(E22-E21)/E21"
$ws.Range("L22").Formula = "=(F22-F21)/F21"

# 7) Fix the "difference vs trend" formulas on row 25 to compare against the
#    new latest-year row (22) instead of the previous year (21).
$ws.Range("H25").Formula = "=H22-H24"
$ws.Range("I25").Formula = "=I22-I24"
$ws.Range("J25").Formula = "=J22-J24"
$ws.Range("K25").Formula = "=K22-K24"
$ws.Range("L25").Formula = "=L22-L24"

# 8) Fix the DCA weighting formulas (row 31) to be driven off the actual new
#    2023 row (22) rather than the extrapolated trend row.
$ws.Range("J31").Formula = "=(K22+L22)/(J22+K22+L22)"
$ws.Range("K31").Formula = "=(J22+L22)/(J22+K22+L22)"
$ws.Range("L31").Formula = "=(J22+K22)/(J22+K22+L22)"

# 8b) Fix the "YTD /1000" helper row (row 28) to point at the new latest-year
#     row (22) instead of the previous year (21), so the CAGR row (29) below
#     it recalculates off the current data.
$ws.Range("H28").Formula = "=(B22/1000)"
$ws.Range("I28").Formula = "=(C22/1000)"
$ws.Range("J28").Formula = "=(D22/1000)"
$ws.Range("K28").Formula = "=(E22/1000)"
$ws.Range("L28").Formula = "=(F22/1000)"

# 9) Update the lump-sum comparison amount.
$ws.Range("C33").Value = 5000

# 10) Restore the last-used selection.
$ws.Range("N19").Select()

$wb.Save()
